$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header (shared string) values
$ws.Range("B1").Value = "TISG_PDR_G"
$ws.Range("D1").Value = "buy_BEE_MWH"
$ws.Range("E1").Value = "sell_lago_MWH"
$ws.Range("F1").Value = "need_to_buy_MW"

# Data rows: row, B, C, D, E, F
$data = @(
    @(2,  6845.33687561269, 5627.73995751213, 9020, 9494.06157,    -30.9806395041899),
    @(3,  1975.76241112899, 3267.83927646187, 9004, 5172.149195,   -105.823914152797),
    @(4,  1903.68386863218, 3361.37557349393, 9004, 5098.786406,   -101.980078714094),
    @(5,  7742.89201065376, 6630.45568318764, 9004, 12164.369865,  85.3305640639116),
    @(6,  7988.16690763696, 7317.60009413061, 9004, 12568.815951,  120.593714062235),
    @(7,  8101.77168568444, 7650.90258619371, 9004, 13070.81297,   150.664327937886),
    @(8,  8220.67689438168, 7945.55096991421, 9004, 13080.859215,  158.405553772189),
    @(9,  8225.40163354142, 7202.17050940992, 9004, 13116.412013,  128.715870369521),
    @(10, 2997.19681285912, 4749.47378762205, 9004, 7400.349805,   6.19278249012207),
    @(11, 2868.31746358633, 4627.32916879301, 9004, 7259.132841,   0.589356091944865),
    @(12, 8947.0418620382,  7745.18551282043, 9004, 14487.69428,   178.409913782593),
    @(13, 8947.04310337663, 8008.9534829775,  9004, 14487.699347,  189.400405275036),
    @(14, 8947.04310337663, 8090.15610275387, 9004, 14487.699347,  192.783847765718),
    @(15, 8947.04310337663, 8053.1902562265,  9004, 14487.699347,  191.243604160411)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
